$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename column B header from "total" to "total_pm10"
$ws.Range("B1").Value = "total_pm10"

# Widen column B to fit the new, longer header text
$ws.Columns.Item(2).ColumnWidth = 16.29

# Page setup (paper size / orientation) as saved by the newer Excel build
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = $excel.xlPortrait

# Leave the active selection on B2, matching the resaved workbook state
[void]$ws.Range("B2").Select()
